# Insert a new data row at row 321 (pushing the existing rows 321-396 down
# to 322-397) and populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 321..396 down to 322..397, leaving a blank row 321 in place.
$ws.Rows.Item(321).Insert()

# Populate the newly inserted row 321 with the new record.
$ws.Cells.Item(321, 1).Value  = 9
$ws.Cells.Item(321, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(321, 3).Value  = "Metropolitana"
$ws.Cells.Item(321, 4).Value  = 45204
$ws.Cells.Item(321, 5).Value  = 13
$ws.Cells.Item(321, 6).Value  = 100112001
$ws.Cells.Item(321, 7).Value  = "Berenjena"
$ws.Cells.Item(321, 8).Value  = "Sin especificar"
$ws.Cells.Item(321, 9).Value  = "Primera"
$ws.Cells.Item(321, 10).Value = 70
$ws.Cells.Item(321, 11).Value = 7000
$ws.Cells.Item(321, 12).Value = 8000
$ws.Cells.Item(321, 13).Value = 7500
$ws.Cells.Item(321, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(321, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(321, 16).Value = 150
$ws.Cells.Item(321, 17).Value = 50
$ws.Cells.Item(321, 18).Value = "Hortaliza"
